$wb = $excel.ActiveWorkbook

# Scheduled market-price refresh: update currentAveragePrice* / Leve*Profit*
# columns (H:N) across all job tables with freshly pulled values.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 208.33333
$ws.Range("H28").Value = 1187.0952
$ws.Range("I28").Value = 959.4375
$ws.Range("J28").Value = 1915.6
$ws.Range("K28").Value = 959.4375
$ws.Range("L28").Value = 1915.6
$ws.Range("M28").Value = -474.4375
$ws.Range("N28").Value = -2885.6
$ws.Range("H29").Value = 2813.6
$ws.Range("I29").Value = 356
$ws.Range("K29").Value = 1068
$ws.Range("M29").Value = -787
$ws.Range("H38").Value = 2142.5334
$ws.Range("I38").Value = 163.14285
$ws.Range("K38").Value = 489.42855
$ws.Range("M38").Value = -117.42855
$ws.Range("H96").Value = 2394.8
$ws.Range("I96").Value = 1991.3334
$ws.Range("K96").Value = 5974.0002
$ws.Range("M96").Value = -4601.0002
$ws.Range("H138").Value = 3258.1667
$ws.Range("I138").Value = 1682.762
$ws.Range("J138").Value = 3993.3555
$ws.Range("K138").Value = 5048.286
$ws.Range("L138").Value = 11980.0665
$ws.Range("M138").Value = 91.71399999999994
$ws.Range("N138").Value = -22260.0665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2327.805
$ws.Range("I32").Value = 1405.2097
$ws.Range("K32").Value = 1405.2097
$ws.Range("M32").Value = -1118.2097
$ws.Range("H45").Value = 2349.8
$ws.Range("I45").Value = 2315.5789
$ws.Range("K45").Value = 2315.5789
$ws.Range("M45").Value = -1938.5789
$ws.Range("H61").Value = 1973.1918
$ws.Range("I61").Value = 1701.7069
$ws.Range("K61").Value = 1701.7069
$ws.Range("M61").Value = -1489.7069
$ws.Range("H74").Value = 3321.8823
$ws.Range("I74").Value = 1997.8334
$ws.Range("J74").Value = 4044.0908
$ws.Range("K74").Value = 1997.8334
$ws.Range("L74").Value = 4044.0908
$ws.Range("M74").Value = -1123.8334
$ws.Range("N74").Value = -5792.0908
$ws.Range("H77").Value = 3321.8823
$ws.Range("I77").Value = 1997.8334
$ws.Range("J77").Value = 4044.0908
$ws.Range("K77").Value = 9989.166999999999
$ws.Range("L77").Value = 20220.454
$ws.Range("M77").Value = -5621.166999999999
$ws.Range("N77").Value = -28956.454
$ws.Range("H132").Value = 2305.7837
$ws.Range("I132").Value = 2111.9412
$ws.Range("K132").Value = 6335.823600000001
$ws.Range("M132").Value = -3805.823600000001
$ws.Range("H136").Value = 1973.1918
$ws.Range("I136").Value = 1701.7069
$ws.Range("K136").Value = 5105.120699999999
$ws.Range("M136").Value = -2555.120699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3945.0952
$ws.Range("I134").Value = 3236.5715
$ws.Range("K134").Value = 9709.7145
$ws.Range("M134").Value = -7174.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 10061.9
$ws.Range("I134").Value = 9985.842000000001
$ws.Range("K134").Value = 29957.526
$ws.Range("M134").Value = -27422.526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 216.66667
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 7500
$ws.Range("L62").Value = 7500
$ws.Range("M62").Value = -6814
$ws.Range("N62").Value = -8872
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 22500
$ws.Range("M65").Value = -19068
$ws.Range("N65").Value = -29364
$ws.Range("H109").Value = 3825.6667
$ws.Range("I109").Value = 3830.8
$ws.Range("K109").Value = 11492.4
$ws.Range("M109").Value = -10452.4
$ws.Range("H112").Value = 2663
$ws.Range("I112").Value = 1327
$ws.Range("K112").Value = 3981
$ws.Range("M112").Value = -2873

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 278.34375
$ws.Range("I2").Value = 198.4
$ws.Range("J2").Value = 411.58334
$ws.Range("K2").Value = 198.4
$ws.Range("L2").Value = 411.58334
$ws.Range("M2").Value = -85.40000000000001
$ws.Range("N2").Value = -637.58334
$ws.Range("H6").Value = 3300
$ws.Range("I6").Value = 3300
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 3300
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -3187
$ws.Range("N6").ClearContents()
$ws.Range("H16").Value = 3300
$ws.Range("I16").Value = 3300
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3300
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3050
$ws.Range("N16").ClearContents()
$ws.Range("H97").Value = 1206.2941
$ws.Range("I97").Value = 1551.1
$ws.Range("K97").Value = 1551.1
$ws.Range("M97").Value = -1055.1
$ws.Range("H132").Value = 17687.346
$ws.Range("I132").Value = 13281.739
$ws.Range("K132").Value = 39845.217
$ws.Range("M132").Value = -37315.217

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1093.7778
$ws.Range("I22").Value = 532.3333
$ws.Range("K22").Value = 532.3333
$ws.Range("M22").Value = -237.3333
$ws.Range("H27").Value = 1093.7778
$ws.Range("I27").Value = 532.3333
$ws.Range("K27").Value = 532.3333
$ws.Range("M27").Value = -425.3333
$ws.Range("H46").Value = 1926.8334
$ws.Range("J46").Value = 2001.2222
$ws.Range("L46").Value = 2001.2222
$ws.Range("N46").Value = -2377.2222
$ws.Range("H61").Value = 4342.1816
$ws.Range("I61").Value = 1774
$ws.Range("J61").Value = 6482.3335
$ws.Range("K61").Value = 1774
$ws.Range("L61").Value = 6482.3335
$ws.Range("M61").Value = -1572
$ws.Range("N61").Value = -6886.3335
$ws.Range("H93").Value = 5599.8
$ws.Range("I93").Value = 5599.8
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 5599.8
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -4351.8
$ws.Range("N93").ClearContents()
$ws.Range("H113").Value = 4342.1816
$ws.Range("I113").Value = 1774
$ws.Range("J113").Value = 6482.3335
$ws.Range("K113").Value = 1774
$ws.Range("L113").Value = 6482.3335
$ws.Range("M113").Value = 396
$ws.Range("N113").Value = -10822.3335
$ws.Range("H132").Value = 2042.2413
$ws.Range("I132").Value = 2042.2413
$ws.Range("K132").Value = 6126.7239
$ws.Range("M132").Value = -3596.7239

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H45").Value = 12284.8
$ws.Range("I45").Value = 8000
$ws.Range("J45").Value = 13356
$ws.Range("K45").Value = 8000
$ws.Range("L45").Value = 13356
$ws.Range("M45").Value = -7509
$ws.Range("N45").Value = -14338
$ws.Range("H122").Value = 2526.5557
$ws.Range("I122").Value = 2105.7144
$ws.Range("J122").Value = 3999.5
$ws.Range("K122").Value = 6317.1432
$ws.Range("L122").Value = 11998.5
$ws.Range("M122").Value = -3867.1432
$ws.Range("N122").Value = -16898.5
$ws.Range("H132").Value = 2899.1428
$ws.Range("I132").Value = 2467.2632
$ws.Range("K132").Value = 7401.7896
$ws.Range("M132").Value = -4871.7896
$ws.Range("H136").Value = 3623.5122
$ws.Range("J136").Value = 9999.333000000001
$ws.Range("L136").Value = 29997.999
$ws.Range("N136").Value = -35097.999
